# Update TPM-derived values on the active sheet (rows 2-4) per new TPM computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.006543333333333334
$ws.Range("H2").Value = 0.01963
$ws.Range("M2").Value = 44.32896033333333
$ws.Range("N2").Value = 132.986881
$ws.Range("O2").Value = 0.5209689208718928
$ws.Range("P2").Value = 0.5209689208718928
$ws.Range("Q2").Value = 0.2900591637811111
$ws.Range("R2").Value = 2.61053247403
$ws.Range("S2").Value = 0.5209689208718928
$ws.Range("T2").Value = 0.5209689208718928

# Row 3
$ws.Range("G3").Value = 0.006543333333333334
$ws.Range("H3").Value = 0.01963
$ws.Range("O3").Value = 0.2245216581053889
$ws.Range("P3").Value = 0.2245216581053888
$ws.Range("Q3").Value = 0.1250066209166667
$ws.Range("R3").Value = 1.12505958825
$ws.Range("S3").Value = 0.2245216581053889
$ws.Range("T3").Value = 0.2245216581053888

# Row 4
$ws.Range("G4").Value = 0.006543333333333334
$ws.Range("H4").Value = 0.01963
$ws.Range("M4").Value = 21.65606733333333
$ws.Range("N4").Value = 64.96820199999999
$ws.Range("O4").Value = 0.2545094210227183
$ws.Range("P4").Value = 0.2545094210227183
$ws.Range("Q4").Value = 0.1417028672511111
$ws.Range("R4").Value = 1.27532580526
$ws.Range("S4").Value = 0.2545094210227183
$ws.Range("T4").Value = 0.2545094210227183
